# Update the "as_of_utc" timestamp column (AA) on the stats sheets
# from "2025-11-21 07:04:04" to "2025-11-21 07:19:58" for every data row.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Главные", "Линейные")
$newTimestamp = "2025-11-21 07:19:58"

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    for ($row = 2; $row -le 26; $row++) {
        $ws.Cells.Item($row, 27).Value = $newTimestamp
    }
}
